$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "season record" header columns (AD:AF), reusing the same
# header formatting (bold, centered, bordered) already applied to A1:AC1
# by copying the format from the adjacent header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2-56) with the
# team's 2017 Oakland Athletics record: 75 wins, 87 losses, 0 ties.
for ($row = 2; $row -le 56; $row++) {
    $ws.Cells.Item($row, 30).Value = 75
    $ws.Cells.Item($row, 31).Value = 87
    $ws.Cells.Item($row, 32).Value = 0
}
